$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.473.70'
$ws.Range('E2').Value = '  +1.52%  '
$ws.Range('D3').Value = '3.795.31'
$ws.Range('E3').Value = '  -0.27%  '
$ws.Range('E4').Value = '  +0.35%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '628.83'
$ws.Range('E5').Value = '  +4.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.43'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('D7').Value = '3.793.39'
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').Value = '  +0.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.452'
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.64'
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  -1.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.60'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').Value = '4.455.87'
$ws.Range('E15').Value = '  +0.42%  '
$ws.Range('D16').Value = '3.852.23'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').Value = '69.666.10'
$ws.Range('E17').Value = '  +2.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.94'
$ws.Range('E18').Value = '  -2.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.11'
$ws.Range('E19').Value = '  +0.02%  '
$ws.Range('E20').Value = '  -0.94%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '467.57'
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.62'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.704'
$ws.Range('E23').Value = '  +0.26%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000148'
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.64'
$ws.Range('E25').Value = '  +0.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.04'
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.15'
$ws.Range('E27').Value = '  +1.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.01'
$ws.Range('E28').Value = '  -0.20%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').Value = '3.966.66'
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.68'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.22'
$ws.Range('E32').Value = '  -0.85%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.28'
$ws.Range('E33').Value = '  -1.56%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '29.00'
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.753.46'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.01'
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.102'
$ws.Range('E38').Value = '  +2.52%  '
$ws.Range('E39').Value = '  +8.51%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.31'
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.86'
$ws.Range('E41').Value = '  +0.46%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.974'
$ws.Range('E42').Value = '  -1.64%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '154.60'
$ws.Range('E45').Value = '  +1.87%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.299'
$ws.Range('E46').Value = '  -0.96%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.94'
$ws.Range('E47').Value = '  +2.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '47.03'
$ws.Range('E48').Value = '  -1.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '42.67'
$ws.Range('E49').Value = '  -1.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.39'
$ws.Range('E50').Value = '  +2.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.42'
$ws.Range('E51').Value = '  +0.40%  '
